$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Minimum" headers to "Minimum Overall"
$ws.Range("J1").Value = "Minimum Overall Supplier Name"
$ws.Range("K1").Value = "Minimum Overall Supplier Unit Price"
$ws.Range("L1").Value = "Minimum Overall Supplier Date"

# Remove the "cabbiage" (row 3) and "onion" (row 3 again, after first shift) rows
# so the rice / white rice / brown rice rows move up to rows 3-5
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()
